$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values in rows 7, 8, 9 (odds data refresh) ---
# Row 7
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 4.75
$ws.Range("I7").Value = 1.5
$ws.Range("J7").Value = 4.75
$ws.Range("L7").Value = 1.91
$ws.Range("U7").Value = 1.5
$ws.Range("V7").Value = 2.5
$ws.Range("W7").Value = 23
$ws.Range("X7").Value = 34
$ws.Range("Y7").Value = 17
$ws.Range("AA7").Value = 34
$ws.Range("AG7").Value = 101
$ws.Range("AK7").Value = 13
$ws.Range("AL7").Value = 11
$ws.Range("AN7").Value = 7.5
$ws.Range("AO7").Value = 23
$ws.Range("AP7").Value = 23
$ws.Range("AQ7").Value = 67
$ws.Range("AU7").Value = 7.5
$ws.Range("AW7").Value = 4
$ws.Range("AX7").Value = 7.5
$ws.Range("AZ7").Value = 19
$ws.Range("BA7").Value = 29

# Row 8
$ws.Range("K8").Value = 2.1
$ws.Range("L8").Value = 3.25
$ws.Range("S8").Value = 1.4
$ws.Range("T8").Value = 2.75
$ws.Range("W8").Value = 8.5
$ws.Range("AK8").Value = 26
$ws.Range("AL8").Value = 21
$ws.Range("AO8").Value = 15
$ws.Range("AT8").Value = 2.75
$ws.Range("AY8").Value = 23

# Row 9
$ws.Range("G9").Value = 1.6
$ws.Range("I9").Value = 5.5
$ws.Range("J9").Value = 2.1
$ws.Range("Q9").Value = 1.57
$ws.Range("R9").Value = 2.35
$ws.Range("U9").Value = 1.62
$ws.Range("V9").Value = 2.2
$ws.Range("X9").Value = 9
$ws.Range("AE9").Value = 13
$ws.Range("AJ9").Value = 17
$ws.Range("AO9").Value = 8
$ws.Range("AU9").Value = 7.5
$ws.Range("AV9").Value = 41
$ws.Range("BA9").Value = 81

# --- Add new rows 11, 12, 13 (additional LaLiga2 matches) ---
# Row 11
$ws.Range("A11").Value = "tjlFtyLO"
$ws.Range("B11").Value = "24/10/2024"
$ws.Range("C11").Value = "16:00"
$ws.Range("D11").Value = "SPAIN - LALIGA2"
$ws.Range("E11").Value = "Burgos CF"
$ws.Range("F11").Value = "Racing Club Ferrol"
$ws.Range("G11").Value = 1.91
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 4.75
$ws.Range("J11").Value = 2.63
$ws.Range("K11").Value = 1.95
$ws.Range("L11").Value = 5
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6.5
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 2.5
$ws.Range("Q11").Value = 2.6
$ws.Range("R11").Value = 1.48
$ws.Range("S11").Value = 1.57
$ws.Range("T11").Value = 2.25
$ws.Range("U11").Value = 2.2
$ws.Range("V11").Value = 1.62
$ws.Range("W11").Value = 5.5
$ws.Range("X11").Value = 7.5
$ws.Range("Y11").Value = 9.5
$ws.Range("Z11").Value = 15
$ws.Range("AA11").Value = 19
$ws.Range("AB11").Value = 41
$ws.Range("AC11").Value = 6
$ws.Range("AD11").Value = 6
$ws.Range("AE11").Value = 19
$ws.Range("AF11").Value = 81
$ws.Range("AG11").Value = 1250
$ws.Range("AH11").Value = 10
$ws.Range("AI11").Value = 23
$ws.Range("AJ11").Value = 17
$ws.Range("AK11").Value = 51
$ws.Range("AL11").Value = 41
$ws.Range("AM11").Value = 51
$ws.Range("AN11").Value = 3.75
$ws.Range("AO11").Value = 11
$ws.Range("AP11").Value = 26
$ws.Range("AQ11").Value = 41
$ws.Range("AR11").Value = 67
$ws.Range("AS11").Value = 251
$ws.Range("AT11").Value = 2.25
$ws.Range("AU11").Value = 9.5
$ws.Range("AV11").Value = 81
$ws.Range("AW11").Value = 6
$ws.Range("AX11").Value = 26
$ws.Range("AY11").Value = 41
$ws.Range("AZ11").Value = 101
$ws.Range("BA11").Value = 151
$ws.Range("BB11").Value = 351
$ws.Range("BC11").Value = 81
$ws.Range("BD11").Value = 81

# Row 12
$ws.Range("A12").Value = "KAGF9ZlB"
$ws.Range("B12").Value = "24/10/2024"
$ws.Range("C12").Value = "16:00"
$ws.Range("D12").Value = "SPAIN - LALIGA2"
$ws.Range("E12").Value = "Gijon"
$ws.Range("F12").Value = "Huesca"
$ws.Range("G12").Value = 1.75
$ws.Range("H12").Value = 3.4
$ws.Range("I12").Value = 4.75
$ws.Range("J12").Value = 2.5
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 5.5
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 2.63
$ws.Range("Q12").Value = 2.35
$ws.Range("R12").Value = 1.57
$ws.Range("S12").Value = 1.5
$ws.Range("T12").Value = 2.5
$ws.Range("U12").Value = 2.2
$ws.Range("V12").Value = 1.62
$ws.Range("W12").Value = 5.5
$ws.Range("X12").Value = 7
$ws.Range("Y12").Value = 9
$ws.Range("Z12").Value = 13
$ws.Range("AA12").Value = 17
$ws.Range("AB12").Value = 34
$ws.Range("AC12").Value = 7
$ws.Range("AD12").Value = 6.5
$ws.Range("AE12").Value = 21
$ws.Range("AF12").Value = 81
$ws.Range("AG12").Value = 201
$ws.Range("AH12").Value = 10
$ws.Range("AI12").Value = 23
$ws.Range("AJ12").Value = 17
$ws.Range("AK12").Value = 51
$ws.Range("AL12").Value = 41
$ws.Range("AM12").Value = 51
$ws.Range("AN12").Value = 3.6
$ws.Range("AO12").Value = 10
$ws.Range("AP12").Value = 26
$ws.Range("AQ12").Value = 34
$ws.Range("AR12").Value = 67
$ws.Range("AS12").Value = 251
$ws.Range("AT12").Value = 2.5
$ws.Range("AU12").Value = 9.5
$ws.Range("AV12").Value = 81
$ws.Range("AW12").Value = 6.5
$ws.Range("AX12").Value = 29
$ws.Range("AY12").Value = 41
$ws.Range("AZ12").Value = 101
$ws.Range("BA12").Value = 151
$ws.Range("BB12").Value = 500
$ws.Range("BC12").Value = 81
$ws.Range("BD12").Value = 81

# Row 13
$ws.Range("A13").Value = "vaTmiRYq"
$ws.Range("B13").Value = "24/10/2024"
$ws.Range("C13").Value = "16:00"
$ws.Range("D13").Value = "SPAIN - LALIGA2"
$ws.Range("E13").Value = "Tenerife"
$ws.Range("F13").Value = "Malaga"
$ws.Range("G13").Value = 2.3
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = 3.4
$ws.Range("J13").Value = 3.2
$ws.Range("K13").Value = 1.91
$ws.Range("L13").Value = 4.33
$ws.Range("M13").Value = 1.13
$ws.Range("N13").Value = 6
$ws.Range("O13").Value = 1.57
$ws.Range("P13").Value = 2.25
$ws.Range("Q13").Value = 2.7
$ws.Range("R13").Value = 1.44
$ws.Range("S13").Value = 1.62
$ws.Range("T13").Value = 2.2
$ws.Range("U13").Value = 2.25
$ws.Range("V13").Value = 1.57
$ws.Range("W13").Value = 5.5
$ws.Range("X13").Value = 9.5
$ws.Range("Y13").Value = 10
$ws.Range("Z13").Value = 21
$ws.Range("AA13").Value = 23
$ws.Range("AB13").Value = 41
$ws.Range("AC13").Value = 6
$ws.Range("AD13").Value = 6
$ws.Range("AE13").Value = 21
$ws.Range("AF13").Value = 81
$ws.Range("AG13").Value = 201
$ws.Range("AH13").Value = 7.5
$ws.Range("AI13").Value = 15
$ws.Range("AJ13").Value = 13
$ws.Range("AK13").Value = 41
$ws.Range("AL13").Value = 34
$ws.Range("AM13").Value = 51
$ws.Range("AN13").Value = 4
$ws.Range("AO13").Value = 15
$ws.Range("AP13").Value = 29
$ws.Range("AQ13").Value = 51
$ws.Range("AR13").Value = 81
$ws.Range("AS13").Value = 301
$ws.Range("AT13").Value = 2.2
$ws.Range("AU13").Value = 9.5
$ws.Range("AV13").Value = 81
$ws.Range("AW13").Value = 5
$ws.Range("AX13").Value = 21
$ws.Range("AY13").Value = 41
$ws.Range("AZ13").Value = 81
$ws.Range("BA13").Value = 126
$ws.Range("BB13").Value = 500
$ws.Range("BC13").Value = 81
$ws.Range("BD13").Value = 81
